# Update Name of Algo - apply new computed values for columns D (Algo column) and E (F column)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -6.952199999999999
$ws.Range("E4").Value = 13.3527

$ws.Range("E5").Value = 13.3199

$ws.Range("D7").Value = -7.512299999999999

$ws.Range("E8").Value = 14.34859999999999

$ws.Range("D16").Value = -8.041099999999998
$ws.Range("E16").Value = 13.32820000000001
